# Update Name of Algo
# Applies corrected RandomForest-imputed values to column B of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B11"  = 6.3695
    "B12"  = 4.956499999999997
    "B15"  = 5.245499999999995
    "B27"  = 6.229900000000002
    "B28"  = 5.561999999999999
    "B31"  = 5.3453
    "B32"  = 6.519299999999996
    "B36"  = 9.425200000000007
    "B38"  = 4.747699999999996
    "B46"  = 6.234600000000001
    "B54"  = 5.033800000000006
    "B55"  = 6.366099999999998
    "B56"  = 4.516299999999994
    "B67"  = 5.653099999999995
    "B69"  = 5.553299999999995
    "B72"  = 5.250000000000004
    "B73"  = 8.9551
    "B83"  = 5.324699999999996
    "B86"  = 5.4256
    "B91"  = 5.032499999999994
    "B93"  = 5.547999999999999
    "B99"  = 5.5066
    "B104" = 9.527599999999996
    "B105" = 8.305300000000006
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
